$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Add new row 5 data ("Car - Germany" sensitivity case)
# -----------------------------------------------------------------
$ws.Range("A5").Value = "None yet"
$ws.Range("B5").Value = "Car - Germany"
$ws.Range("C5").Value = "Ford Fusion"
$ws.Range("D5").Value = "Ford Fusion Hybrid"
$ws.Range("E5").Value = "mile"
$ws.Range("F5").Value = "gal"
$ws.Range("G5").Value = 126.62163000000001
$ws.Range("H5").Value = 3.3893390630606466
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 6.1834509000000004
$ws.Range("K5").Value = -0.2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 25
$ws.Range("O5").Value = 42
$ws.Range("P5").Value = 14425
$ws.Range("Q5").Value = 27401.277693029457
$ws.Range("R5").Value = 28216.1
$ws.Range("S5").Value = 7
$ws.Range("T5").Value = 27523.4
$ws.Range("U5").Value = 7
$ws.Range("V5").Value = 2861.13425677328
$ws.Range("W5").Value = 2774.6681207909455
$ws.Range("X5").Value = 34000
$ws.Range("Y5").Value = 14
$ws.Range("Z5").Value = 40000
$ws.Range("AA5").Value = 14

# -----------------------------------------------------------------
# Add new row 6 data ("Lamp - Germany" sensitivity case)
# -----------------------------------------------------------------
$ws.Range("A6").Value = "None yet"
$ws.Range("B6").Value = "Lamp - Germany"
$ws.Range("C6").Value = "Incandescent"
$ws.Range("D6").Value = "LED"
$ws.Range("E6").Value = "lm-hr"
$ws.Range("F6").Value = "kW-hr"
$ws.Range("G6").Value = 3.6
$ws.Range("H6").Value = 3.3893390630606466
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0.36299999999999999
$ws.Range("K6").Value = -0.4
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 8833.3333333333303
$ws.Range("O6").Value = 81800
$ws.Range("P6").Value = 580350
$ws.Range("Q6").Value = 27401.277693029457
$ws.Range("R6").Value = 1.88
$ws.Range("S6").Value = 1.8
$ws.Range("T6").Value = 1.21
$ws.Range("U6").Value = 10
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 2.2000000000000002
$ws.Range("Y6").Value = 1.8
$ws.Range("Z6").Value = 6.5
$ws.Range("AA6").Value = 10

# -----------------------------------------------------------------
# Comments (units) for row 5, mirroring the existing row 3 ("Car") comments
# -----------------------------------------------------------------
$ws.Range("G5").AddComment("Matthew Heun:`nMJ/gallon") | Out-Null
$ws.Range("H5").AddComment("Matthew Heun:`nMJ/$") | Out-Null
$ws.Range("J5").AddComment("Matthew Heun:`n$/gal") | Out-Null
$ws.Range("N5").AddComment("Matthew Heun:`nmi/gallon") | Out-Null
$ws.Range("O5").AddComment("Matthew Heun:`nmi/gallon") | Out-Null
$ws.Range("P5").AddComment("Matthew Heun:`nmi/yr") | Out-Null
$ws.Range("Q5").AddComment("Matthew Heun:`n$/year") | Out-Null
$ws.Range("R5").AddComment("Matthew Heun:`n$") | Out-Null
$ws.Range("S5").AddComment("Matthew Heun:`nyear") | Out-Null
$ws.Range("T5").AddComment("Matthew Heun:`n$") | Out-Null
$ws.Range("U5").AddComment("Matthew Heun:`nyear") | Out-Null
$ws.Range("V5").AddComment("Matthew Heun:`n$/year") | Out-Null
$ws.Range("W5").AddComment("Matthew Heun:`n$/year") | Out-Null
$ws.Range("X5").AddComment("Matthew Heun:`nMJ") | Out-Null
$ws.Range("Y5").AddComment("Matthew Heun:`nyears") | Out-Null
$ws.Range("Z5").AddComment("Matthew Heun:`nMJ") | Out-Null
$ws.Range("AA5").AddComment("Matthew Heun:`nyears") | Out-Null

# -----------------------------------------------------------------
# Comments (units) for row 6, mirroring the existing row 4 ("Lamp") comments
# -----------------------------------------------------------------
$ws.Range("G6").AddComment("Matthew Heun:`nMJ/kW-hr") | Out-Null
$ws.Range("H6").AddComment("Matthew Heun:`nMJ/$") | Out-Null
$ws.Range("J6").AddComment("Matthew Heun:`n$/kW-hr") | Out-Null
$ws.Range("N6").AddComment("Matthew Heun:`nlm-hr/kW-hr") | Out-Null
$ws.Range("O6").AddComment("Matthew Heun:`nlm-hr/kW-hr") | Out-Null
$ws.Range("P6").AddComment("Matthew Heun:`nLm-hr/yr") | Out-Null
$ws.Range("Q6").AddComment("Matthew Heun:`n$/year") | Out-Null
$ws.Range("R6").AddComment("Matthew Heun:`n$") | Out-Null
$ws.Range("S6").AddComment("Matthew Heun:`nyear") | Out-Null
$ws.Range("T6").AddComment("Matthew Heun:`n$") | Out-Null
$ws.Range("U6").AddComment("Matthew Heun:`nyear") | Out-Null
$ws.Range("V6").AddComment("Matthew Heun:`n$/year") | Out-Null
$ws.Range("W6").AddComment("Matthew Heun:`n$/year") | Out-Null
$ws.Range("X6").AddComment("Matthew Heun:`nMJ") | Out-Null
$ws.Range("Y6").AddComment("Matthew Heun:`nyears") | Out-Null
$ws.Range("Z6").AddComment("Matthew Heun:`nMJ") | Out-Null
$ws.Range("AA6").AddComment("Matthew Heun:`nyears") | Out-Null

# -----------------------------------------------------------------
# Restore the T5 direct-formatting style used by sibling row 3 (T3),
# matching the 11pt font cell style already used elsewhere in the sheet.
# -----------------------------------------------------------------
$ws.Range("T5").Font.Size = 11

# -----------------------------------------------------------------
# Final selection, matching the saved workbook's active cell (J5)
# -----------------------------------------------------------------
$ws.Range("J5").Select()
